$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1272.4
$ws.Range("I15").Value = 1272.4
$ws.Range("K15").Value = 3817.2
$ws.Range("M15").Value = -3648.2
$ws.Range("H40").Value = 2793.3635
$ws.Range("I40").Value = 2099.5
$ws.Range("J40").Value = 2947.5557
$ws.Range("K40").Value = 2099.5
$ws.Range("L40").Value = 2947.5557
$ws.Range("M40").Value = -1924.5
$ws.Range("N40").Value = -3297.5557
$ws.Range("H70").Value = 2898.7
$ws.Range("J70").Value = 3548.3635
$ws.Range("L70").Value = 10645.0905
$ws.Range("N70").Value = -11185.0905
$ws.Range("H73").Value = 2898.7
$ws.Range("J73").Value = 3548.3635
$ws.Range("L73").Value = 10645.0905
$ws.Range("N73").Value = -12517.0905
$ws.Range("H103").Value = 739.8
$ws.Range("J103").Value = 999.6667
$ws.Range("L103").Value = 2999.0001
$ws.Range("N103").Value = -4171.0001
$ws.Range("H111").Value = 1375
$ws.Range("I111").Value = 1375
$ws.Range("K111").Value = 4125
$ws.Range("M111").Value = -1058
$ws.Range("H125").Value = 10474.833
$ws.Range("I125").Value = 9712.5
$ws.Range("K125").Value = 87412.5
$ws.Range("M125").Value = -84952.5
$ws.Range("H137").Value = 1886.6666
$ws.Range("I137").Value = 1477.8
$ws.Range("J137").Value = 3931
$ws.Range("K137").Value = 4433.4
$ws.Range("L137").Value = 11793
$ws.Range("M137").Value = -1883.4
$ws.Range("N137").Value = -16893

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8244.172
$ws.Range("I32").Value = 7016.0586
$ws.Range("K32").Value = 7016.0586
$ws.Range("M32").Value = -6729.0586
$ws.Range("H74").Value = 1014.8
$ws.Range("I74").Value = 976.9167
$ws.Range("J74").Value = 1166.3334
$ws.Range("K74").Value = 976.9167
$ws.Range("L74").Value = 1166.3334
$ws.Range("M74").Value = -102.9167
$ws.Range("N74").Value = -2914.3334
$ws.Range("H77").Value = 1014.8
$ws.Range("I77").Value = 976.9167
$ws.Range("J77").Value = 1166.3334
$ws.Range("K77").Value = 4884.5835
$ws.Range("L77").Value = 5831.666999999999
$ws.Range("M77").Value = -516.5834999999997
$ws.Range("N77").Value = -14567.667
$ws.Range("H110").Value = 3131.524
$ws.Range("I110").Value = 1707
$ws.Range("K110").Value = 1707
$ws.Range("M110").Value = 338

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3725.5
$ws.Range("I105").Value = 3678.6667
$ws.Range("J105").Value = 3785.7144
$ws.Range("K105").Value = 3678.6667
$ws.Range("L105").Value = 3785.7144
$ws.Range("M105").Value = -1931.6667
$ws.Range("N105").Value = -7279.7144
$ws.Range("H134").Value = 9249.75
$ws.Range("I134").Value = 9249.75
$ws.Range("K134").Value = 27749.25
$ws.Range("M134").Value = -25214.25

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 172.2381
$ws.Range("I2").Value = 184.92857
$ws.Range("K2").Value = 1109.57142
$ws.Range("M2").Value = -996.57142
$ws.Range("H109").Value = 2780
$ws.Range("I109").Value = 892
$ws.Range("J109").Value = 7500
$ws.Range("K109").Value = 2676
$ws.Range("L109").Value = 22500
$ws.Range("M109").Value = -1636
$ws.Range("N109").Value = -24580
$ws.Range("H138").Value = 1407.4166
$ws.Range("I138").Value = 599.25
$ws.Range("J138").Value = 1811.5
$ws.Range("K138").Value = 1797.75
$ws.Range("L138").Value = 5434.5
$ws.Range("M138").Value = 3342.25
$ws.Range("N138").Value = -15714.5

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 15818.182
$ws.Range("I46").Value = 6333.3335
$ws.Range("J46").Value = 19375
$ws.Range("K46").Value = 6333.3335
$ws.Range("L46").Value = 19375
$ws.Range("M46").Value = -6177.3335
$ws.Range("N46").Value = -19687
$ws.Range("H107").Value = 265.8
$ws.Range("I107").Value = 274.75
$ws.Range("K107").Value = 274.75
$ws.Range("M107").Value = 1645.25
$ws.Range("H126").Value = 1132.8
$ws.Range("I126").Value = 1137.5
$ws.Range("K126").Value = 3412.5
$ws.Range("M126").Value = -942.5
$ws.Range("H140").Value = 60000
$ws.Range("J140").Value = 60000
$ws.Range("L140").Value = 60000
$ws.Range("N140").Value = -70360
$ws.Range("H141").Value = 34999.5
$ws.Range("J141").Value = 34999.5
$ws.Range("L141").Value = 34999.5
$ws.Range("N141").Value = -45359.5

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1385.409
$ws.Range("I46").Value = 832.4167
$ws.Range("J46").Value = 2049
$ws.Range("K46").Value = 832.4167
$ws.Range("L46").Value = 2049
$ws.Range("M46").Value = -644.4167
$ws.Range("N46").Value = -2425
$ws.Range("H132").Value = 14210.964
$ws.Range("I132").Value = 16819.47
$ws.Range("J132").Value = 10179.637
$ws.Range("K132").Value = 50458.41
$ws.Range("L132").Value = 30538.911
$ws.Range("M132").Value = -47928.41
$ws.Range("N132").Value = -35598.911
$ws.Range("H136").Value = 3612.9167
$ws.Range("I136").Value = 3435
$ws.Range("J136").Value = 4502.5
$ws.Range("K136").Value = 10305
$ws.Range("L136").Value = 13507.5
$ws.Range("M136").Value = -7755
$ws.Range("N136").Value = -18607.5

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 2000
$ws.Range("J26").Value = 2000
$ws.Range("L26").Value = 2000
$ws.Range("N26").Value = -2586
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H107").Value = 137.5
$ws.Range("I107").Value = 125
$ws.Range("K107").Value = 375
$ws.Range("M107").Value = 1545
$ws.Range("H113").Value = 2254.182
$ws.Range("I113").Value = 828.7143
$ws.Range("J113").Value = 4748.75
$ws.Range("K113").Value = 2486.1429
$ws.Range("L113").Value = 14246.25
$ws.Range("M113").Value = -316.1428999999998
$ws.Range("N113").Value = -18586.25
$ws.Range("H135").Value = 671124.25
$ws.Range("I135").Value = 2509499.5
$ws.Range("K135").Value = 2509499.5
$ws.Range("M135").Value = -2504429.5
$ws.Range("H136").Value = 3636.7334
$ws.Range("I136").Value = 3636.7334
$ws.Range("K136").Value = 10910.2002
$ws.Range("M136").Value = -8360.200199999999
$ws.Range("H140").Value = 68666.336
$ws.Range("J140").Value = 68666.336
$ws.Range("L140").Value = 68666.336
$ws.Range("N140").Value = -79026.336
